# 996: Merge ND0 files
# Insert a new "ARMS" worksheet immediately before "T2A", populate its header
# row and formatting, and keep "T2A" as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new ARMS sheet right before T2A -------------------------
$t2a = $wb.Worksheets.Item("T2A")
$arms = $wb.Worksheets.Add($t2a)
$arms.Name = "ARMS"

# --- 2. Header row values ----------------------------------------------------
$headers = @(
    "Assessment_Staff_Name",
    "Assessment_Staff_Key",
    "Assessment_Staff_Grade",
    "Assessmentent_Team_Key",
    "Assessment_Provider_Code",
    "CRN",
    "Disposal_or_Release_Date",
    "Sentence_Type",
    "SO_Registration_Date"
)
$cols = @("A","B","C","D","E","F","G","H","I")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $arms.Range("$($cols[$i])1").Value = $headers[$i]
}

$arms.Range("A1:I1").RowHeight = 37

# --- 3. Formatting ------------------------------------------------------------
# Colours (COM colors are BGR-encoded longs)
$blue  = 10909496   # FF3877A6
$gray  = 11642277   # FFA5A5B1
$fill  = 10511371   # FF0B64A0
$white = 16777215   # FFFFFFFF

$headerRange = $arms.Range("A1:I1")
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 9
$headerRange.Font.Bold = $true
$headerRange.Font.Color = $white
$headerRange.Interior.Color = $fill
$headerRange.Interior.PatternColor = $white
$headerRange.HorizontalAlignment = -4131   # xlLeft
$headerRange.WrapText = $true

# Column C ("Assessment_Staff_Grade") is stored as text and is not wrapped
$arms.Range("C1").NumberFormat = "@"
$arms.Range("C1").WrapText = $false

# Borders -- apply per cell (Excel treats Left/Right on a multi-cell range as
# the outer edge of the whole range, not each cell, so loop individually).
# A1:E1 -> thin blue top/right, thin grey bottom; A1 also gets a thin blue left
# F1:I1 -> thin blue right only
for ($i = 0; $i -lt 5; $i++) {
    $cell = $arms.Range("$($cols[$i])1")
    $cell.Borders.Item(10).Color = $blue    # xlEdgeRight
    $cell.Borders.Item(8).Color = $blue     # xlEdgeTop
    $cell.Borders.Item(9).Color = $gray     # xlEdgeBottom
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
}

$arms.Range("A1").Borders.Item(7).Color = $blue   # xlEdgeLeft
$arms.Range("A1").Borders.Item(7).LineStyle = 1

for ($i = 5; $i -lt 9; $i++) {
    $cell = $arms.Range("$($cols[$i])1")
    $cell.Borders.Item(10).Color = $blue
    $cell.Borders.Item(10).LineStyle = 1
}

# --- 4. Keep T2A as the selected / active sheet -------------------------------
$t2aFresh = $wb.Worksheets.Item("T2A")
$t2aFresh.Activate()
$t2aFresh.Range("A1:AO1").Select()
